$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Metadata sheet updates
$ws1.Range("B3").Value = "2.0.0"
$ws1.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$ws1.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# Include from SNOMED CT sheet updates
# Update the first concept code (kept as text, matching the existing column's
# string-typed cells), then remove the second concept row entirely
$scratch = $ws2.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "116224001"
$scratch.Copy()
$ws2.Range("A2").PasteSpecial(-4163)
$scratch.Clear()

$ws2.Rows(3).Delete()
